# Applies the scheduled-runner profit-sheet refresh described in the commit diff.
# For each affected row (one "flip" entry per worksheet), this rewrites the
# buy/sell-price and profit columns (H:N) with the refreshed market values.

$wb = $excel.ActiveWorkbook

# --- ALC sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 1002930.1
$ws.Range("I43").Value = 2250.25
$ws.Range("J43").Value = 1670050
$ws.Range("K43").Value = 2250.25
$ws.Range("L43").Value = 1670050
$ws.Range("M43").Value = -2181.25
$ws.Range("N43").Value = -1670188
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 5620738.5
$ws.Range("J112").Value = 5955232.5
$ws.Range("L112").Value = 17865697.5
$ws.Range("N112").Value = -17867913.5
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 125005210
$ws.Range("I137").Value = 100006330
$ws.Range("J137").Value = 166670020
$ws.Range("K137").Value = 300018990
$ws.Range("L137").Value = 500010060
$ws.Range("M137").Value = -300016440
$ws.Range("N137").Value = -500015160
# Row 138: All-night Crafting
$ws.Range("H138").Value = 3975342.8
$ws.Range("J138").Value = 4451854
$ws.Range("L138").Value = 13355562
$ws.Range("N138").Value = -13365842

# --- ARM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 15881516
$ws.Range("I32").Value = 16956008
$ws.Range("K32").Value = 16956008
$ws.Range("M32").Value = -16955721
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1932.7778
$ws.Range("I45").Value = 1420.6666
$ws.Range("J45").Value = 2188.8333
$ws.Range("K45").Value = 1420.6666
$ws.Range("L45").Value = 2188.8333
$ws.Range("M45").Value = -1043.6666
$ws.Range("N45").Value = -2942.8333
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 1833.3334
$ws.Range("I63").Value = 1833.3334
$ws.Range("K63").Value = 1833.3334
$ws.Range("M63").Value = -1147.3334
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 1833.3334
$ws.Range("I66").Value = 1833.3334
$ws.Range("K66").Value = 9166.666999999999
$ws.Range("M66").Value = -5734.666999999999
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 2298.4546
$ws.Range("I102").Value = 2129.2
$ws.Range("J102").Value = 3991
$ws.Range("K102").Value = 2129.2
$ws.Range("L102").Value = 3991
$ws.Range("M102").Value = -507.1999999999998
$ws.Range("N102").Value = -7235
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 25607.688
$ws.Range("I110").Value = 26981.533
$ws.Range("K110").Value = 26981.533
$ws.Range("M110").Value = -24936.533
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2863.2415
$ws.Range("I122").Value = 2405.7144
$ws.Range("K122").Value = 7217.1432
$ws.Range("M122").Value = -4767.1432
# Row 133: Shielding My Students
$ws.Range("H133").Value = 85271.414
$ws.Range("I133").Value = 64134.5
$ws.Range("J133").Value = 89498.8
$ws.Range("K133").Value = 64134.5
$ws.Range("L133").Value = 89498.8
$ws.Range("M133").Value = -61604.5
$ws.Range("N133").Value = -94558.8
# Row 134: Brace for More Vambraces
$ws.Range("H134").Value = 445000
$ws.Range("J134").Value = 445000
$ws.Range("L134").Value = 445000
$ws.Range("N134").Value = -455140
# Row 139: Backing up My Words
$ws.Range("H139").Value = 68856.89
$ws.Range("J139").Value = 68856.89
$ws.Range("L139").Value = 68856.89
$ws.Range("N139").Value = -79136.89

# --- BSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Row 60: Talon Terrors
$ws.Range("H60").Value = 74999.5
$ws.Range("I60").Value = 70000
$ws.Range("J60").Value = 79999
$ws.Range("K60").Value = 70000
$ws.Range("L60").Value = 79999
$ws.Range("M60").Value = -69401
$ws.Range("N60").Value = -81197
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 3573
$ws.Range("I99").Value = 1475
$ws.Range("J99").Value = 6370.3335
$ws.Range("K99").Value = 1475
$ws.Range("L99").Value = 6370.3335
$ws.Range("M99").Value = 23
$ws.Range("N99").Value = -9366.333500000001
# Row 138: Bladewinner
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

# --- CRP sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 2214.923
$ws.Range("I16").Value = 2340.0833
$ws.Range("J16").Value = 713
$ws.Range("K16").Value = 2340.0833
$ws.Range("L16").Value = 713
$ws.Range("M16").Value = -2053.0833
$ws.Range("N16").Value = -1287
# Row 31: Wall Not Found
$ws.Range("H31").Value = 31254824
$ws.Range("I31").Value = 3935.48
$ws.Range("J31").Value = 142865140
$ws.Range("K31").Value = 3935.48
$ws.Range("L31").Value = 142865140
$ws.Range("M31").Value = -3640.48
$ws.Range("N31").Value = -142865730
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 31254824
$ws.Range("I34").Value = 3935.48
$ws.Range("J34").Value = 142865140
$ws.Range("K34").Value = 3935.48
$ws.Range("L34").Value = 142865140
$ws.Range("M34").Value = -3733.48
$ws.Range("N34").Value = -142865544
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3797.7222
$ws.Range("I58").Value = 4056.8125
$ws.Range("J58").Value = 1725
$ws.Range("K58").Value = 4056.8125
$ws.Range("L58").Value = 1725
$ws.Range("M58").Value = -3853.8125
$ws.Range("N58").Value = -2131
# Row 69: Landing the Big One
$ws.Range("H69").Value = 165999.8
$ws.Range("I69").Value = 119999.664
$ws.Range("K69").Value = 119999.664
$ws.Range("M69").Value = -119250.664
# Row 72: Fishing for Profits (L)
$ws.Range("H72").Value = 165999.8
$ws.Range("I72").Value = 119999.664
$ws.Range("K72").Value = 359998.992
$ws.Range("M72").Value = -356254.992
# Row 107: Built to Last
$ws.Range("H107").Value = 2692
$ws.Range("I107").Value = 2247
$ws.Range("K107").Value = 2247
$ws.Range("M107").Value = -327
# Row 113: Patient Patients
$ws.Range("H113").Value = 2214.923
$ws.Range("I113").Value = 2340.0833
$ws.Range("J113").Value = 713
$ws.Range("K113").Value = 2340.0833
$ws.Range("L113").Value = 713
$ws.Range("M113").Value = -170.0832999999998
$ws.Range("N113").Value = -5053
# Row 136: Turali Quality
$ws.Range("H136").Value = 3797.7222
$ws.Range("I136").Value = 4056.8125
$ws.Range("J136").Value = 1725
$ws.Range("K136").Value = 12170.4375
$ws.Range("L136").Value = 5175
$ws.Range("M136").Value = -9620.4375
$ws.Range("N136").Value = -10275
# Row 138: Bow Out
$ws.Range("H138").Value = 88193.5
$ws.Range("J138").Value = 88193.5
$ws.Range("L138").Value = 88193.5
$ws.Range("N138").Value = -98473.5
# Row 140: Spear Pressure
$ws.Range("H140").Value = 67500
# Row 141: No Greater Treasure
$ws.Range("H141").Value = 439803.56
$ws.Range("J141").Value = 458615.53
$ws.Range("L141").Value = 458615.53
$ws.Range("N141").Value = -468975.53

# --- CUL sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 1934.1428
$ws.Range("I3").Value = 923.1667
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 2769.5001
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = -2657.5001
$ws.Range("N3").Value = -24224
# Row 123: Topping Up the Pot
$ws.Range("H123").Value = 2660
$ws.Range("I123").Value = 2660
$ws.Range("K123").Value = 7980
$ws.Range("M123").Value = -5530
# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 7091.9653
$ws.Range("I134").Value = 1508.45
$ws.Range("K134").Value = 4525.35
$ws.Range("M134").Value = 544.6499999999996
# Row 137: Creative Chocolate
$ws.Range("H137").Value = 5052.1333
$ws.Range("I137").Value = 2528.3
$ws.Range("J137").Value = 10099.8
$ws.Range("K137").Value = 7584.900000000001
$ws.Range("L137").Value = 30299.4
$ws.Range("M137").Value = -2484.900000000001
$ws.Range("N137").Value = -40499.39999999999

# --- GSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 5000
$ws.Range("I97").Value = 4000
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 4000
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -3504
$ws.Range("N97").Value = -6992
# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 87997.625
$ws.Range("J135").Value = 87997.625
$ws.Range("L135").Value = 87997.625
$ws.Range("N135").Value = -98137.625

# --- LTW sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 61995
$ws.Range("I74").Value = 61995
$ws.Range("K74").Value = 61995
$ws.Range("M74").Value = -60997
# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 61995
$ws.Range("I77").Value = 61995
$ws.Range("K77").Value = 185985
$ws.Range("M77").Value = -180993
# Row 131: For What Was Gleaned
$ws.Range("H131").Value = 53993.668
$ws.Range("J131").Value = 53993.668
$ws.Range("L131").Value = 53993.668
$ws.Range("N131").Value = -64073.668
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 10666
$ws.Range("I136").Value = 9336
$ws.Range("J136").Value = 11996
$ws.Range("K136").Value = 28008
$ws.Range("L136").Value = 35988
$ws.Range("M136").Value = -25458
$ws.Range("N136").Value = -41088
